# The upstream data source renamed the "ScreenRecStarted" reading-state
# category to "0_unstated". Propagate that rename to every cell (header +
# transition-matrix row labels) that still carries the old name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains("ScreenRecStarted")) {
            $cell.Value2 = $val.Replace("ScreenRecStarted", "0_unstated")
        }
    }
}

$ws.Range("A27").Select()
